# Insert a new data row before the current row 28 (shifts old rows 28-64 down
# to 29-65) and populate the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this pushes the existing row 28 (and all
# rows below it) down by one, copying formatting from the row above/below as
# Excel normally does for a row insert.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record's values.
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 44587
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112040
$ws.Cells.Item(28, 7).Value = "Cilantro"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 300
$ws.Cells.Item(28, 11).Value = 2300
$ws.Cells.Item(28, 12).Value = 2500
$ws.Cells.Item(28, 13).Value = 2400
$ws.Cells.Item(28, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 1200
$ws.Cells.Item(28, 17).Value = 2
$ws.Cells.Item(28, 18).Value = "Hortaliza"
